$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new data row at the bottom (A60) with value "GRT-USD"
$ws.Range("A60").Value = "GRT-USD"
